$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two rows that were dropped ("05-00-00" / 팀빌딩 워크숍, and "09-00-00" / 회의록 작성) ---
# Delete from the bottom up so row numbers above the deleted row stay valid.
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(20).Delete()

# After the two deletions the sheet now reads (rows 2-22):
#  20 = 06-00-00, 21 = 07-00-00, 22 = 08-00-00

# --- Insert two new sub-task rows under "07-00-00" ---
$ws.Rows.Item(22).Insert()
$ws.Rows.Item(22).Insert()

$ws.Range("A22").Value = "07-01-00"
$ws.Range("B22").Value = "메일 송부"
$ws.Range("C22").Value = " "
$ws.Range("D22").Value = " "
$ws.Range("E22").Value = "진행"
$ws.Range("F22").Value = " "
$ws.Range("G22").Value = " "
$ws.Range("H22").Value = " "
$ws.Range("I22").Value = 2

$ws.Range("A23").Value = "07-01-01"
$ws.Range("B23").Value = " "
$ws.Range("C23").Value = " "
$ws.Range("D23").Value = " "
$ws.Range("E23").Value = "진행"
$ws.Range("F23").Value = " "
$ws.Range("G23").Value = " "
$ws.Range("H23").Value = " "
$ws.Range("I23").Value = 3

# "08-00-00" is now on row 24 (unchanged content, just shifted up by the earlier deletes)

# --- Insert nine new sub-task rows under "08-00-00" (new team rotation rows) ---
for ($i = 0; $i -lt 9; $i++) {
    $ws.Rows.Item(25).Insert()
}

$ws.Range("A25").Value = "08-01-00"
$ws.Range("B25").Value = "완료"
$ws.Range("C25").Value = "차체설계1팀"
$ws.Range("D25").Value = " "
$ws.Range("E25").Value = "진행"
$ws.Range("F25").Value = " "
$ws.Range("G25").Value = " "
$ws.Range("H25").Value = " "
$ws.Range("I25").Value = 2

$ws.Range("A26").Value = "08-02-00"
$ws.Range("B26").Value = " "
$ws.Range("C26").Value = "차체설계2팀"
$ws.Range("D26").Value = " "
$ws.Range("E26").Value = "진행"
$ws.Range("F26").Value = " "
$ws.Range("G26").Value = " "
$ws.Range("H26").Value = " "
$ws.Range("I26").Value = 2

$ws.Range("A27").Value = "08-03-00"
$ws.Range("B27").Value = " "
$ws.Range("C27").Value = "차체설계3팀"
$ws.Range("D27").Value = " "
$ws.Range("E27").Value = "진행"
$ws.Range("F27").Value = " "
$ws.Range("G27").Value = " "
$ws.Range("H27").Value = " "
$ws.Range("I27").Value = 2

$ws.Range("A28").Value = "08-04-00"
$ws.Range("B28").Value = " "
$ws.Range("C28").Value = "외장설계1팀"
$ws.Range("D28").Value = " "
$ws.Range("E28").Value = "진행"
$ws.Range("F28").Value = " "
$ws.Range("G28").Value = " "
$ws.Range("H28").Value = " "
$ws.Range("I28").Value = 2

$ws.Range("A29").Value = "08-05-00"
$ws.Range("B29").Value = " "
$ws.Range("C29").Value = "dkkkkkkkkkkkkkkkkk"
$ws.Range("D29").Value = " "
$ws.Range("E29").Value = "진행"
$ws.Range("F29").Value = " "
$ws.Range("G29").Value = " "
$ws.Range("H29").Value = " "
$ws.Range("I29").Value = 2

$ws.Range("A30").Value = "08-06-00"
$ws.Range("B30").Value = " "
$ws.Range("C30").Value = "daaleiw12222222222222222222222"
$ws.Range("D30").Value = " "
$ws.Range("E30").Value = "진행"
$ws.Range("F30").Value = " "
$ws.Range("G30").Value = " "
$ws.Range("H30").Value = " "
$ws.Range("I30").Value = 2

$ws.Range("A31").Value = "08-07-00"
$ws.Range("B31").Value = " "
$ws.Range("C31").Value = "aksdakalskdfasf"
$ws.Range("D31").Value = " "
$ws.Range("E31").Value = "진행"
$ws.Range("F31").Value = " "
$ws.Range("G31").Value = " "
$ws.Range("H31").Value = " "
$ws.Range("I31").Value = 2

$ws.Range("A32").Value = "08-08-00"
$ws.Range("B32").Value = " "
$ws.Range("C32").Value = "1212123124k"
$ws.Range("D32").Value = " "
$ws.Range("E32").Value = "진행"
$ws.Range("F32").Value = " "
$ws.Range("G32").Value = " "
$ws.Range("H32").Value = " "
$ws.Range("I32").Value = 2

$ws.Range("A33").Value = "08-09-00"
$ws.Range("B33").Value = " "
$ws.Range("C33").Value = "외장설계2팀"
$ws.Range("D33").Value = " "
$ws.Range("E33").Value = "진행"
$ws.Range("F33").Value = " "
$ws.Range("G33").Value = " "
$ws.Range("H33").Value = " "
$ws.Range("I33").Value = 2
